$wb = $excel.ActiveWorkbook

# --- ALC (sheet index 1) ---
$ws = $wb.Worksheets.Item(1)
# row 63
$ws.Range("H63").Value = 20000
$ws.Range("I63").Value = 20000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -19376
$ws.Range("N63").ClearContents()
# row 66
$ws.Range("H66").Value = 20000
$ws.Range("I66").Value = 20000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 60000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -56880
$ws.Range("N66").ClearContents()
# row 137
$ws.Range("H137").Value = 26317674
$ws.Range("I137").Value = 1229.4
$ws.Range("J137").Value = 125004340
$ws.Range("K137").Value = 3688.2
$ws.Range("L137").Value = 375013020
$ws.Range("M137").Value = -1138.2
$ws.Range("N137").Value = -375018120

# --- ARM (sheet index 2) ---
$ws = $wb.Worksheets.Item(2)
# row 2
$ws.Range("H2").Value = 892070.9
$ws.Range("I2").Value = 647.7826
$ws.Range("J2").Value = 2942344
$ws.Range("K2").Value = 647.7826
$ws.Range("L2").Value = 2942344
$ws.Range("M2").Value = -534.7826
$ws.Range("N2").Value = -2942570
# row 32
$ws.Range("H32").Value = 4635.103
$ws.Range("I32").Value = 4748.7856
$ws.Range("K32").Value = 4748.7856
$ws.Range("M32").Value = -4461.7856
# row 45
$ws.Range("H45").Value = 1045.1538
$ws.Range("I45").Value = 878.55554
$ws.Range("J45").Value = 1420
$ws.Range("K45").Value = 878.55554
$ws.Range("L45").Value = 1420
$ws.Range("M45").Value = -501.55554
$ws.Range("N45").Value = -2174
# row 74
$ws.Range("H74").Value = 8727.833000000001
$ws.Range("I74").Value = 1257.4286
$ws.Range("J74").Value = 13481.728
$ws.Range("K74").Value = 1257.4286
$ws.Range("L74").Value = 13481.728
$ws.Range("M74").Value = -383.4286
$ws.Range("N74").Value = -15229.728
# row 77
$ws.Range("H77").Value = 8727.833000000001
$ws.Range("I77").Value = 1257.4286
$ws.Range("J77").Value = 13481.728
$ws.Range("K77").Value = 6287.143
$ws.Range("L77").Value = 67408.64
$ws.Range("M77").Value = -1919.143
$ws.Range("N77").Value = -76144.64
# row 116
$ws.Range("H116").Value = 892070.9
$ws.Range("I116").Value = 647.7826
$ws.Range("J116").Value = 2942344
$ws.Range("K116").Value = 647.7826
$ws.Range("L116").Value = 2942344
$ws.Range("M116").Value = 1646.2174
$ws.Range("N116").Value = -2946932

# --- BSM (sheet index 3) ---
$ws = $wb.Worksheets.Item(3)
# row 3
$ws.Range("H3").Value = 892070.9
$ws.Range("I3").Value = 647.7826
$ws.Range("J3").Value = 2942344
$ws.Range("K3").Value = 647.7826
$ws.Range("L3").Value = 2942344
$ws.Range("M3").Value = -533.7826
$ws.Range("N3").Value = -2942572
# row 99
$ws.Range("H99").Value = 1059.6666
$ws.Range("I99").Value = 1051.8
$ws.Range("J99").Value = 1069.5
$ws.Range("K99").Value = 1051.8
$ws.Range("L99").Value = 1069.5
$ws.Range("M99").Value = 446.2
$ws.Range("N99").Value = -4065.5

# --- CRP (sheet index 4) ---
$ws = $wb.Worksheets.Item(4)
# row 31
$ws.Range("H31").Value = 1479.8518
$ws.Range("I31").Value = 1129.2632
$ws.Range("J31").Value = 2312.5
$ws.Range("K31").Value = 1129.2632
$ws.Range("L31").Value = 2312.5
$ws.Range("M31").Value = -834.2632000000001
$ws.Range("N31").Value = -2902.5
# row 34
$ws.Range("H34").Value = 1479.8518
$ws.Range("I34").Value = 1129.2632
$ws.Range("J34").Value = 2312.5
$ws.Range("K34").Value = 1129.2632
$ws.Range("L34").Value = 2312.5
$ws.Range("M34").Value = -927.2632000000001
$ws.Range("N34").Value = -2716.5

# --- CUL (sheet index 5) ---
$ws = $wb.Worksheets.Item(5)
# row 62
$ws.Range("H62").Value = 3666.6667
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3666.6667
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11000.0001
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -12372.0001
# row 65
$ws.Range("H65").Value = 3666.6667
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3666.6667
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 33000.0003
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -39864.0003
# row 98
$ws.Range("H98").Value = 900
$ws.Range("I98").Value = 1150
$ws.Range("J98").Value = 400
$ws.Range("K98").Value = 3450
$ws.Range("L98").Value = 1200
$ws.Range("M98").Value = -1952
$ws.Range("N98").Value = -4196
# row 131
$ws.Range("H131").Value = 3186.4285
$ws.Range("I131").Value = 4869.1816
$ws.Range("J131").Value = 2589.3225
$ws.Range("K131").Value = 14607.5448
$ws.Range("L131").Value = 7767.967500000001
$ws.Range("M131").Value = -9567.5448
$ws.Range("N131").Value = -17847.9675

# --- GSM (sheet index 6) ---
$ws = $wb.Worksheets.Item(6)
# row 102
$ws.Range("H102").Value = 1156.2
$ws.Range("I102").Value = 1061.6923
$ws.Range("J102").Value = 1331.7142
$ws.Range("K102").Value = 1061.6923
$ws.Range("L102").Value = 1331.7142
$ws.Range("M102").Value = 560.3077000000001
$ws.Range("N102").Value = -4575.7142

# --- LTW (sheet index 7) ---
$ws = $wb.Worksheets.Item(7)
# row 22
$ws.Range("H22").Value = 510.7143
$ws.Range("I22").Value = 465.83334
$ws.Range("J22").Value = 544.375
$ws.Range("K22").Value = 465.83334
$ws.Range("L22").Value = 544.375
$ws.Range("M22").Value = -170.83334
$ws.Range("N22").Value = -1134.375
# row 27
$ws.Range("H27").Value = 510.7143
$ws.Range("I27").Value = 465.83334
$ws.Range("J27").Value = 544.375
$ws.Range("K27").Value = 465.83334
$ws.Range("L27").Value = 544.375
$ws.Range("M27").Value = -358.83334
$ws.Range("N27").Value = -758.375
# row 61
$ws.Range("H61").Value = 1141.6428
$ws.Range("I61").Value = 1051.8422
$ws.Range("J61").Value = 1331.2222
$ws.Range("K61").Value = 1051.8422
$ws.Range("L61").Value = 1331.2222
$ws.Range("M61").Value = -849.8422
$ws.Range("N61").Value = -1735.2222
# row 113
$ws.Range("H113").Value = 1141.6428
$ws.Range("I113").Value = 1051.8422
$ws.Range("J113").Value = 1331.2222
$ws.Range("K113").Value = 1051.8422
$ws.Range("L113").Value = 1331.2222
$ws.Range("M113").Value = 1118.1578
$ws.Range("N113").Value = -5671.2222

# --- WVR (sheet index 8) ---
$ws = $wb.Worksheets.Item(8)
# row 133
$ws.Range("H133").Value = 55357.5
$ws.Range("J133").Value = 55357.5
$ws.Range("L133").Value = 55357.5
$ws.Range("N133").Value = -65477.5
# row 140
$ws.Range("H140").Value = 17776.334
$ws.Range("J140").Value = 17776.334
$ws.Range("L140").Value = 17776.334
$ws.Range("N140").Value = -28136.334
